# Refresh the crypto symbol price/volume snapshot (GitHub Actions bot run).
# All Price (D) / Volume(1h) (E) cells hold plain text (e.g. "307.57",
# "-0.12%"), not numbers. Assigning a numeric-looking string straight to
# .Value lets Excel auto-coerce it to a real number/percentage, which would
# silently change the stored representation. Forcing the cell to Text format
# before the assignment keeps it a literal string, and ClearFormats()
# afterwards drops the temporary "@" number format again so the cell's style
# ends up back at its original (default) style instead of picking up a new one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "307.57"
Set-TextValue $ws.Range("E2") "-0.12%"
Set-TextValue $ws.Range("D3") "41.06"
Set-TextValue $ws.Range("E3") "-0.05%"
Set-TextValue $ws.Range("D4") "5.250"
Set-TextValue $ws.Range("E4") "2.31%"
Set-TextValue $ws.Range("D5") "0.07666"
Set-TextValue $ws.Range("E5") "0.47%"
Set-TextValue $ws.Range("D6") "1.626"
Set-TextValue $ws.Range("E6") "0.32%"
Set-TextValue $ws.Range("D7") "0.9181"
Set-TextValue $ws.Range("E7") "1.86%"
Set-TextValue $ws.Range("D9") "0.1236"
Set-TextValue $ws.Range("E9") "13.00%"
Set-TextValue $ws.Range("D10") "0.1828"
Set-TextValue $ws.Range("E10") "2.78%"
Set-TextValue $ws.Range("D11") "0.09114"
Set-TextValue $ws.Range("E11") "-0.57%"
Set-TextValue $ws.Range("D12") "0.04272"
Set-TextValue $ws.Range("E12") "1.76%"
Set-TextValue $ws.Range("D13") "0.1051"
Set-TextValue $ws.Range("E13") "0.00%"
Set-TextValue $ws.Range("D14") "0.001262"
Set-TextValue $ws.Range("E14") "1.02%"
Set-TextValue $ws.Range("D15") "0.005811"
Set-TextValue $ws.Range("E15") "0.18%"
Set-TextValue $ws.Range("D18") "4.317"
Set-TextValue $ws.Range("E18") "1.25%"
Set-TextValue $ws.Range("D19") "0.3335"
Set-TextValue $ws.Range("D20") "7.402"
Set-TextValue $ws.Range("E20") "13.28%"
Set-TextValue $ws.Range("D21") "0.1383"
Set-TextValue $ws.Range("E21") "1.73%"
Set-TextValue $ws.Range("D22") "0.2893"
Set-TextValue $ws.Range("E22") "2.88%"
Set-TextValue $ws.Range("D23") "0.04069"
Set-TextValue $ws.Range("E23") "-0.51%"
Set-TextValue $ws.Range("D25") "0.004356"
Set-TextValue $ws.Range("E25") "8.90%"
Set-TextValue $ws.Range("D26") "0.0001273"
Set-TextValue $ws.Range("E26") "-2.06%"
Set-TextValue $ws.Range("D38") "0.02471"
Set-TextValue $ws.Range("E38") "3.51%"
Set-TextValue $ws.Range("D39") "0.05282"
Set-TextValue $ws.Range("E39") "1.94%"
Set-TextValue $ws.Range("D40") "0.007848"
Set-TextValue $ws.Range("E40") "0.86%"
Set-TextValue $ws.Range("E41") "1.08%"
Set-TextValue $ws.Range("D42") "0.006506"
Set-TextValue $ws.Range("E42") "-6.58%"
Set-TextValue $ws.Range("D43") "0.001915"
Set-TextValue $ws.Range("E43") "-1.80%"
Set-TextValue $ws.Range("D44") "0.007725"
Set-TextValue $ws.Range("E44") "-9.62%"
Set-TextValue $ws.Range("D45") "0.3059"
Set-TextValue $ws.Range("E45") "-0.46%"
Set-TextValue $ws.Range("D46") "0.00006745"
Set-TextValue $ws.Range("E46") "-2.33%"
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.21%"
Set-TextValue $ws.Range("D48") "0.1699"
Set-TextValue $ws.Range("E48") "850.44%"
Set-TextValue $ws.Range("E49") "-2.43%"
Set-TextValue $ws.Range("D50") "0.00002105"
Set-TextValue $ws.Range("E50") "0.21%"
Set-TextValue $ws.Range("D51") "0.0002005"
Set-TextValue $ws.Range("E51") "0.21%"
